# Populate the previously-empty last row of the "Asset Log" table with the
# jQuery asset entry (Asset / Source / Reasoning / Date Retrieved).

$d = $word.ActiveDocument
$table = $d.Tables.Item(1)

$lastRow = $table.Rows.Count

$table.Cell($lastRow, 1).Range.Text = "js/jQuery.min.js"
$table.Cell($lastRow, 2).Range.Text = "https://jquery.com"
$table.Cell($lastRow, 3).Range.Text = "For document manipulation, traversal and event handling"
$table.Cell($lastRow, 4).Range.Text = "09/01/2024"
